$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# QA test-case sheet reshuffle:
#   - Drop the "US11TimeOut is commented out / long wait" note (row 16, col B)
#     since that comment no longer applies.
#   - Expand the placeholder rows that used to just say "US13" / "US15" /
#     "US17" / "US19" / "US21" / "US23" / "US25" (rows 18-30) into concrete,
#     named test cases for adding/deleting Drugs, Physicians and Patients,
#     each marked PASS, pushing the table out to row 35.
# ---------------------------------------------------------------------------

# Row 16: remove the stale "Testing this can be done..." comment.
$ws.Cells.Item(16, 2).ClearContents()

# Rows 18-35: new/renumbered User Story test cases (column A = test case
# name, column C = result). None of these rows carry a comment.
$newRows = @(
    @{ Row = 18; Name = "US13AddDrugSingle" },
    @{ Row = 19; Name = "US13AddDrugBulk" },
    @{ Row = 20; Name = "US14ViewDrugs" },
    @{ Row = 21; Name = "US15DeleteDrugSingle" },
    @{ Row = 22; Name = "US15DeleteDrugMulti" },
    @{ Row = 23; Name = "US16EditDrug" },
    @{ Row = 24; Name = "US17AddPhysicianSingle" },
    @{ Row = 25; Name = "US17AddPhysicianBulk" },
    @{ Row = 26; Name = "US18ViewPhysicians" },
    @{ Row = 27; Name = "US19DeletePhysicianSingle" },
    @{ Row = 28; Name = "US19DeletePhysicianMulti" },
    @{ Row = 29; Name = "US20EditPhysician" },
    @{ Row = 30; Name = "US21AddPatientSingle" },
    @{ Row = 31; Name = "US21AddPatientBulk" },
    @{ Row = 32; Name = "US22ViewPatients" },
    @{ Row = 33; Name = "US23DeletePatientSingle" },
    @{ Row = 34; Name = "US23DeletePatientMulti" },
    @{ Row = 35; Name = "US24EditPatient" }
)

foreach ($item in $newRows) {
    $ws.Cells.Item($item.Row, 1).Value = $item.Name
    $ws.Cells.Item($item.Row, 3).Value = "PASS"
}

# The green/red PASS-FAIL conditional formatting was manually dragged down
# over a big range; keep it extended to cover the freshly added rows.
$cf = $ws.Range("C2:C631").FormatConditions
for ($i = 1; $i -le $cf.Count; $i++) {
    $null = $cf.Item($i).ModifyAppliesToRange($ws.Range("C2:C637"))
}

# Scroll back to the top and land the selection on the last edited cell.
$null = $ws.Range("B33").Select()
